# Insert a new row for "Sharing half of global tax with low-income countries"
# directly below the existing "Global tax on millionaires" row (row 4),
# shifting the remaining rows down by one. The new row's data values are
# left as placeholders (0) since the figures haven't been computed yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 4 (and everything below it) down by one row.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with its label and placeholder values.
$ws.Range("A4").Value = "Sharing half of global tax with low-income countries"
$ws.Range("B4:G4").Value = 0
